# Apply updates to Ref_LDV ZEV QC workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ZEV market share cap for 2035-2050 from 99% to 100%
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 1
$ws.Range("V3").Value = 1
$ws.Range("W3").Value = 1

# Update the comment text in X3 to reflect the revised wording
$ws.Range("X3").Value = "Should move 2015 & 2020 values (actual sales) once include subsidy policy! ; https://www150.statcan.gc.ca/t1/tbl1/en/cv.action?pid=2010002101"
